# "creating a list of ticker symbols" -- the Ticker Symbol column (A) is
# removed entirely and the 2013 EPS values (column B) shift left into
# column A. The first (most negative) value is highlighted with a red
# fill to call it out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A ("Ticker Symbol" header + AAL, ARNC, ... tickers).
# Column B ("2013" header + the negative EPS values) shifts into column A.
$ws.Columns("A").Delete()

# Flag the largest-magnitude negative EPS value with a solid red fill.
$ws.Range("A2").Interior.Color = 255
